$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)

# Add a paragraph border (all four edges, 5-twip space, no line) to paragraph 1.
$b = $p1.Format.Borders(1)
$b.DistanceFromTop = 5
$b.DistanceFromLeft = 5
$b.DistanceFromBottom = 5
$b.DistanceFromRight = 5

# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25

# Replace the ID placeholder text and drop the trailing space run in one go.
$d.Content.Find.Execute("**ID__AFFARS_5325_topic_6__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5325_204__ID**", 2)
